$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 40, shifting existing rows 40-41 down to 41-42.
$ws.Rows.Item(40).Insert()

# New row 40 data (weekly update)
$ws.Cells.Item(40, 1).Value = 9
$ws.Cells.Item(40, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(40, 3).Value = "Metropolitana"
$ws.Cells.Item(40, 4).Value = 44747
$ws.Cells.Item(40, 4).NumberFormat = $ws.Cells.Item(41, 4).NumberFormat
$ws.Cells.Item(40, 5).Value = 13
$ws.Cells.Item(40, 6).Value = 100112035
$ws.Cells.Item(40, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 43
$ws.Cells.Item(40, 11).Value = 22000
$ws.Cells.Item(40, 12).Value = 22000
$ws.Cells.Item(40, 13).Value = 22000
$ws.Cells.Item(40, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(40, 15).Value = "Hijuelas"
$ws.Cells.Item(40, 16).Value = 1467
$ws.Cells.Item(40, 17).Value = 15
$ws.Cells.Item(40, 18).Value = "Hortaliza"

$ws.Range("A1").Select()
